$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.711.64'
$ws.Range("E2").Value = '  +1.86%  '
$ws.Range("D3").Value = '1.637.67'
$ws.Range("E3").Value = '  +2.10%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("E6").Value = '  +2.12%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.253'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0625'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0840'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.26%  '
$ws.Range("D12").Value = '1.865.28'
$ws.Range("E12").Value = '  +2.06%  '
$ws.Range("D13").Value = '1.622.19'
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("E14").Value = '  +1.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.528'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.03%  '
$ws.Range("D16").Value = '26.734.14'
$ws.Range("E16").Value = '  +2.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.78%  '
$ws.Range("D18").Value = '0.0₃0742'
$ws.Range("E18").Value = '  +2.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '209.75'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.48%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.77%  '
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.121'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0522'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.42%  '
$ws.Range("E32").Value = '  +1.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.97'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("E34").Value = '  +1.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("D36").Value = '1.170.94'
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0168'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.811'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.40%  '
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.506'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.93%  '
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("E42").Value = '  +1.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.01%  '
$ws.Range("D44").Value = '1.777.00'
$ws.Range("E44").Value = '  +2.13%  '
$ws.Range("E45").Value = '  +0.68%  '
$ws.Range("E46").Value = '  +1.54%  '
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.40%  '
$ws.Range("E49").Value = '  +1.68%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.60%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.409'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.61%  '
